$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph right before the table: replace <w:lang w:val="en-US"/>
#    with <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/>
#    NOTE: this must run BEFORE any $d.Tables access - once a Table
#    object has been touched, $d.Paragraphs indices before the table
#    become unreliable in this runtime.
# ---------------------------------------------------------------------
$preP = $d.Paragraphs.Item(32)
$preXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="WW-"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr></w:p>'
$preP.Range.InsertXML($preXml)

# ---------------------------------------------------------------------
# 2) Table column widths (affects gridCol + every tcW in that column)
# ---------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$widths = @(2268, 1843, 709, 1417, 1560, 1417, 1418, 1842, 2552)
for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
    $tbl.Columns.Item($c).Width = $widths[$c - 1] / 20.0
}

# ---------------------------------------------------------------------
# 3) Paragraph alignment fixes in the 3 data rows (rows 2..4):
#    - columns 1,2,4,5,6,7,8,9 : jc="both" removed (-> left/default)
#    - column 3 (DAYS)          : jc="both" -> jc="center"
# ---------------------------------------------------------------------
for ($r = 2; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        if ($c -eq 3) {
            $cell.Range.Paragraphs.Item(1).Alignment = 1
        } else {
            $cell.Range.Paragraphs.Item(1).Alignment = 0
        }
    }
}

# ---------------------------------------------------------------------
# 4) Text content changes in the 2nd row (the placeholder row):
#    ${POSITION_ORG}   -> ${SERVICE_ORG}   (split across 3 runs)
#    ${POSITION_SERVE} -> ${SERVICE_SERVE}/${POSITION} (split across 4 runs)
# ---------------------------------------------------------------------
$cell8 = $tbl.Cell(2, 8)
$xml8 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="WW-"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t>SERVICE</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>_ORG}</w:t></w:r>' +
    '</w:p>'
$cell8.Range.InsertXML($xml8)

$cell9 = $tbl.Cell(2, 9)
$xml9 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="WW-"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t>SERVICE</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>_SERVE}</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t>/${POSITION}</w:t></w:r>' +
    '</w:p>'
$cell9.Range.InsertXML($xml9)
